$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.947.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.560.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.07"
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0598"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.783.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.560.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.970.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.71%  "
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("E33").Value = "  +3.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.418.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("E36").Value = "  +9.75%  "
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.533"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.808"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.696.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0957"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.06%  "
